$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Clean up the redundant "no border" explicit style on B3,C3,B6,C6,B7,C7,
#    B9,C9 -- these cells had an explicit-but-empty border applied; resetting
#    them to the Normal style removes the redundant style definition while
#    leaving their values/visual appearance unchanged.
# ---------------------------------------------------------------------------
foreach ($addr in @("B3","C3","B6","C6","B7","C7","B9","C9")) {
    $ws.Range($addr).Style = "Normal"
}

# ---------------------------------------------------------------------------
# 2) New model-run rows (12-15 and 18), added below the original table.
# ---------------------------------------------------------------------------

# Row 12
$ws.Range("A12").Formula = "=EXP(15)"
$ws.Range("B12").Value = 0.27
$ws.Range("C12").Formula = "=EXP(13.7)"
$ws.Range("D12").Formula = "=B12/(1+(B12+C12)/A12)"
$ws.Range("E12").Value = "if C is similar to N_eggs or slightly larger, then  survival is normal ish. "

# Row 13
$ws.Range("A13").Formula = "=EXP(15)"
$ws.Range("B13").Value = 0.27
$ws.Range("C13").Formula = "=EXP(20)"
$ws.Range("E13").Value = "if n_eggs > C, survival gets really small"

# Row 14
$ws.Range("A14").Formula = "=EXP(15)"
$ws.Range("B14").Value = 0.27
$ws.Range("C14").Formula = "=EXP(14)"

# D13:D14 share one formula group
$ws.Range("D13:D14").Formula = "=B13/(1+(B13+C13)/A13)"

# Row 15
$ws.Range("A15").Formula = "=EXP(15)"
$ws.Range("B15").Value = 0.27
$ws.Range("C15").Formula = "=EXP(25)"
$ws.Range("D15").Formula = "=B15/(1+(B15+C15)/A15)"

# Row 18
$ws.Range("A18").Formula = "=EXP(14)"
$ws.Range("B18").Value = 0.75
$ws.Range("C18").Formula = "=EXP(15)"
$ws.Range("D18").Formula = "=B18/(1+(B18+C18)/A18)"

# ---------------------------------------------------------------------------
# 3) Re-apply the recurring "box" border look to column D of the new rows by
#    copying formats from the existing rows that already carry it -- this
#    reuses the existing style/border table instead of creating new ones.
#    D3/D6/D7/D9 -> right border only ("mid-block")
#    D4/D8       -> right+bottom border ("end-of-block")
# ---------------------------------------------------------------------------
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null

$ws.Range("D4").Copy() | Out-Null
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Column widths for the new, wider data (A/C best-fit-ish, D and E fixed).
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 11.330729166666666
$ws.Columns("C").ColumnWidth = 11.330729166666666
$ws.Columns("D").ColumnWidth = 19.166666666666668
$ws.Columns("E").ColumnWidth = 35.166666666666664

# ---------------------------------------------------------------------------
# 5) Sheet view: zoomed to 150% and the selection moved to A19 (just below
#    the new data), matching where the user was last working.
# ---------------------------------------------------------------------------
$ws.Range("A19").Select() | Out-Null
$excel.ActiveWindow.Zoom = 150
